$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: find the (empty) paragraph that immediately follows the
# "Introduction:" Heading1 paragraph, and rebuild it into two paragraphs of
# new body text, moving the lone "_GoBack" bookmark so it ends up at the end
# of the second new paragraph.
# ---------------------------------------------------------------------------
$introHeadingIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Introduction*") {
        $introHeadingIndex = $i
        break
    }
}

if ($introHeadingIndex -gt 0) {
    $targetIndex = $introHeadingIndex + 1
    $targetPara = $d.Paragraphs.Item($targetIndex)
    $targetRange = $targetPara.Range

    $newXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:eastAsia="zh-CN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:eastAsia="zh-CN"/>
              </w:rPr>
              <w:t xml:space="preserve">In Assignment 3 Where’s the File, we implement another version of git by using network I/O and multithread programing. </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:eastAsia="zh-CN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:eastAsia="zh-CN"/>
              </w:rPr>
              <w:t xml:space="preserve">First, we successfully </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:eastAsia="zh-CN"/>
              </w:rPr>
              <w:t xml:space="preserve">fully implemented this project, which we believe is remarkable and </w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

    $targetRange.InsertXML($newXml)
}

# ---------------------------------------------------------------------------
# Step 2: the document previously carried a single "_GoBack" bookmark inside
# an otherwise-empty paragraph near the end (right before the "Data Section
# specification..." paragraph). Since that bookmark has now been recreated
# above, strip it out of its old location, leaving a plain empty paragraph.
# Search from the end of the document so the freshly-inserted bookmark
# (near the top) is not mistaken for it.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$oldBookmarkIndex = -1
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.WordOpenXML -like "*_GoBack*") {
        $oldBookmarkIndex = $i
        break
    }
}

if ($oldBookmarkIndex -gt 0) {
    $oldPara = $d.Paragraphs.Item($oldBookmarkIndex)
    $oldRange = $oldPara.Range

    $emptyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

    $oldRange.InsertXML($emptyXml)
}
